$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: the "Atas/Bandy/Bradin/.../Wright" author list is re-saved with an extra
# whitespace character inserted between each name (elastic-search re-export
# reformatted the padding), producing a new distinct shared string.
$ws.Range("E2").Value = '[Atas%Jenny%coreGivesNoEmail%1,                          Bandy%Kenneth%coreGivesNoEmail%1,                          Bradin%Stuart A.%coreGivesNoEmail%1,                          Cadwallender%Bruce A.%coreGivesNoEmail%1,                          Cinti%Sandro K.%coreGivesNoEmail%1,                          Collins%Curtis D.%coreGivesNoEmail%1,                          Goldberg%Janet%coreGivesNoEmail%1,                          Holmes%Jennifer G.%coreGivesNoEmail%1,                          Kim%Christopher%coreGivesNoEmail%1,                          Krupansky%Frank%coreGivesNoEmail%1,                          Lozon%Marie M.%coreGivesNoEmail%1,                          Rodgers%Phillip E.%coreGivesNoEmail%1,                          Shlafer%Jean%coreGivesNoEmail%1,                          Wagner%Deborah%coreGivesNoEmail%1,                          Wilkerson%William M.%coreGivesNoEmail%1,                          Wright%Carrie M.%coreGivesNoEmail%1]'

# E4: the "Cathy Campbell / Marianne Baernholdt" author list picks up the same
# kind of re-padding, matching an already-present shared-string variant.
$ws.Range("E4").Value = '[Cathy%Campbell%xref no email%1,       Marianne%Baernholdt%xref no email%1]'
